# Actualización automática 2025-11-24 12:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": nueva venta de INODOROS y PORCELANATO para
#     TULCAN NARVAEZ EDITH MARITZA (fila 22), con su conteo de asesores (fila 23) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H22").Value = 2735.25
$wsGrupo.Range("M22").Value = 12050.73
$wsGrupo.Range("H23").Value = "1 de 21"
$wsGrupo.Range("M23").Value = "4 de 21"

# --- Sheet "VENTA MENSUAL": actualiza la venta de noviembre para el mismo
#     asesor y el total de la columna ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F22").Value = 13481.94
$wsMensual.Range("F23").Value = 36254.99

# --- Sheet "CUMPLIMIENTO MENSUAL": refleja la nueva venta en INODOROS
#     (fila 6), PORCELANATO (fila 12) y el TOTAL (fila 14) ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D6").Value = 2735.25
$wsCumplimiento.Range("E6").Value = 172.3336814602599
$wsCumplimiento.Range("F6").Value = 0.9407295884348513

$wsCumplimiento.Range("D12").Value = 23060.36
$wsCumplimiento.Range("E12").Value = 21357.64
$wsCumplimiento.Range("F12").Value = 0.5191670043675988

$wsCumplimiento.Range("D14").Value = 36254.99000000001
$wsCumplimiento.Range("E14").Value = 19144.48101170094
$wsCumplimiento.Range("F14").Value = 0.6544284509926562

# Ensanchar ligeramente la columna D para que quepa el nuevo valor
# (el offset de 5/6 compensa la conversion interna de Excel entre el ancho
# de columna en caracteres y el ancho de archivo OOXML)
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.166666666666666
